$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 225; existing rows 225-300 shift down to 226-301.
$ws.Rows.Item(225).Insert()

# Populate the newly inserted row 225 with the new data record.
$ws.Range("A225").Value = 11
$ws.Range("B225").Value = "Vega Monumental Concepción"
$ws.Range("C225").Value = "Bíobío"
$ws.Range("D225").Value = 44784
$ws.Range("E225").Value = 8
$ws.Range("F225").Value = 100114001
$ws.Range("G225").Value = "Papa"
$ws.Range("H225").Value = "Asterix"
$ws.Range("I225").Value = "1a (guarda lavada)"
$ws.Range("J225").Value = 4000
$ws.Range("K225").Value = 8500
$ws.Range("L225").Value = 9000
$ws.Range("M225").Value = 8750
$ws.Range("N225").Value = "$/malla 25 kilos"
$ws.Range("O225").Value = "Región de La Araucanía"
$ws.Range("P225").Value = 350
$ws.Range("Q225").Value = 25
$ws.Range("R225").Value = "Hortaliza"
